$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("errName4to50", "Name can only has 4 to 50 charactors.", "ชื่อมีความยาวตั้งแต่ 4 ถึง 50 ตัวอักษร", "Name can only has 4 to 50 charactors.", "Name can only has 4 to 50 charactors.", "Name can only has 4 to 50 charactors.", "Name can only has 4 to 50 charactors.", "Name can only has 4 to 50 charactors.", "Name can only has 4 to 50 charactors."),
    @("errInvalidEmail", "Please enter a valid email address.", "รูปแบบอีเมล์ไม่ถูกต้อง", "Please enter a valid email address.", "Please enter a valid email address.", "Please enter a valid email address.", "Please enter a valid email address.", "Please enter a valid email address.", "Please enter a valid email address."),
    @("errInvalidPhone", "Please enter a valid phone number.", "รูปแบบเบอร์ติดต่อไม่ถูกต้อง", "Please enter a valid phone number.", "Please enter a valid phone number.", "Please enter a valid phone number.", "Please enter a valid phone number.", "Please enter a valid phone number.", "Please enter a valid phone number."),
    @("errAddress20to400", "Address must be between 20 and 400 characters", "ที่อยู่สั้นหรือยาวเกินไป", "Address must be between 20 and 400 characters", "Address must be between 20 and 400 characters", "Address must be between 20 and 400 characters", "Address must be between 20 and 400 characters", "Address must be between 20 and 400 characters", "Address must be between 20 and 400 characters"),
    @("remark", "Remark", "บันทึก", "Remark", "Remark", "Remark", "Remark", "Remark", "Remark"),
    @("remarkPH", "Message to the seller", "ข้อความถึงผู้ขาย", "Message to the seller", "Message to the seller", "Message to the seller", "Message to the seller", "Message to the seller", "Message to the seller"),
    @("error", "Error", "ผิดพลาด", "Error", "Error", "Error", "Error", "Error", "Error"),
    @("errTest", "Test Error", "ทดสอบผิดพลาด", "Test Error", "Test Error", "Test Error", "Test Error", "Test Error", "Test Error"),
    @("errMissing", "Some information is missing", "ข้อมูลบางส่วนสูญหาย", "Some information is missing", "Some information is missing", "Some information is missing", "Some information is missing", "Some information is missing", "Some information is missing"),
    @("errInputMissing", "Some information is missing (input)", "ข้อมูลบางส่วนสูญหาย (input)", "Some information is missing (input)", "Some information is missing (input)", "Some information is missing (input)", "Some information is missing (input)", "Some information is missing (input)", "Some information is missing (input)"),
    @("errCartMissing", "Some information is missing (cart)", "ข้อมูลบางส่วนสูญหาย (cart)", "Some information is missing (cart)", "Some information is missing (cart)", "Some information is missing (cart)", "Some information is missing (cart)", "Some information is missing (cart)", "Some information is missing (cart)"),
    @("errInvalidProdectOptId", "errInvalidProdectOptId", "errInvalidProdectOptId", "errInvalidProdectOptId", "errInvalidProdectOptId", "errInvalidProdectOptId", "errInvalidProdectOptId", "errInvalidProdectOptId", "errInvalidProdectOptId"),
    @("errTotalAmtMismatch", "Total Amount Mismatch.", "ผลรวมไม่สมเหตุผล", "Total Amount Mismatch.", "Total Amount Mismatch.", "Total Amount Mismatch.", "Total Amount Mismatch.", "Total Amount Mismatch.", "Total Amount Mismatch."),
    @("errGrandTotalMismatch", "Grand Total Mismatch", "ผลรวมสุดท้ายไม่สมเหตุผล", "Grand Total Mismatch", "Grand Total Mismatch", "Grand Total Mismatch", "Grand Total Mismatch", "Grand Total Mismatch", "Grand Total Mismatch"),
    @("useQRPay", "Switch to QR pay", "เปลี่ยนเป็นโอนด้วยคิวอาร์", "Switch to QR pay", "Switch to QR pay", "Switch to QR pay", "Switch to QR pay", "Switch to QR pay", "Switch to QR pay"),
    @("useManualPay", "Switch to manual pay", "เปลี่ยนเป็นโอนด้วยเลขบัญชี", "Switch to manual pay", "Switch to manual pay", "Switch to manual pay", "Switch to manual pay", "Switch to manual pay", "Switch to manual pay"),
    @("clickToUpload", "Click to upload", "คลิ๊กเพื่อส่งรูป", "Click to upload", "Click to upload", "Click to upload", "Click to upload", "Click to upload", "Click to upload"),
    @("totalForPay", "Total Pay", "ยอดรวม", "Total Pay", "Total Pay", "Total Pay", "Total Pay", "Total Pay", "Total Pay"),
)

$startRow = 47
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newRows[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}

$win = $excel.ActiveWindow
$win.ScrollRow = 48
$win.ScrollColumn = 1
$ws.Range("D65").Select() | Out-Null
